# Generate Report for Handback
# Update the Correspond Handoff / Handback DateTime values on the
# per-language sheets ("zh-cn" and "de-de") to reflect the latest
# report generation run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 08:59:10"
$wsZhCn.Range("H2").Value = "2016-03-21 08:59:30"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 08:59:13"
$wsDeDe.Range("H2").Value = "2016-03-21 08:59:36"
